$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 2 4 "68.298.45"
Set-TextValue 2 5 "  +8.70%  "
Set-TextValue 3 4 "3.634.75"
Set-TextValue 3 5 "  +4.73%  "
Set-TextValue 4 4 "0.999"
Set-TextValue 4 5 "  -0.10%  "
Set-TextValue 5 4 "420.74"
Set-TextValue 5 5 "  +1.38%  "
Set-TextValue 6 4 "131.86"
Set-TextValue 6 5 "  +0.92%  "
Set-TextValue 7 4 "0.654"
Set-TextValue 7 5 "  +4.05%  "
Set-TextValue 8 4 "3.626.44"
Set-TextValue 8 5 "  +4.71%  "
Set-TextValue 9 5 "  -0.04%  "
Set-TextValue 10 4 "0.778"
Set-TextValue 10 5 "  +6.74%  "
Set-TextValue 11 4 "0.182"
Set-TextValue 11 5 "  +20.10%  "
Set-TextValue 12 4 "0.0000359"
Set-TextValue 12 5 "  +61.59%  "
Set-TextValue 13 4 "42.93"
Set-TextValue 13 5 "  +0.60%  "
Set-TextValue 14 4 "10.04"
Set-TextValue 14 5 "  +3.13%  "
Set-TextValue 15 4 "4.167.37"
Set-TextValue 15 5 "  +3.49%  "
Set-TextValue 16 5 "  -0.13%  "
Set-TextValue 17 4 "20.59"
Set-TextValue 17 5 "  +0.20%  "
Set-TextValue 18 4 "3.649.32"
Set-TextValue 18 5 "  +4.96%  "
Set-TextValue 19 5 "  +5.10%  "
Set-TextValue 20 4 "68.065.37"
Set-TextValue 20 5 "  +8.35%  "
Set-TextValue 21 4 "12.49"
Set-TextValue 21 5 "  -1.31%  "
Set-TextValue 22 4 "467.61"
Set-TextValue 22 5 "  -1.52%  "
Set-TextValue 23 4 "89.58"
Set-TextValue 23 5 "  -1.27%  "
Set-TextValue 24 4 "3.15"
Set-TextValue 24 5 "  -4.44%  "
Set-TextValue 25 4 "13.46"
Set-TextValue 25 5 "  +2.07%  "
Set-TextValue 26 2 "PancakeSwap"
Set-TextValue 26 3 "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue 26 4 "3.39"
Set-TextValue 26 5 "  +1.53%  "
Set-TextValue 27 2 "Filecoin"
Set-TextValue 27 3 "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue 27 4 "10.27"
Set-TextValue 27 5 "  -2.63%  "
Set-TextValue 28 5 "  +8.88%  "
Set-TextValue 29 4 "4.89"
Set-TextValue 29 5 "  +2.02%  "
Set-TextValue 30 5 "  +4.12%  "
Set-TextValue 31 5 "  +2.44%  "
Set-TextValue 32 4 "7.46"
Set-TextValue 32 5 "  -1.44%  "
Set-TextValue 33 5 "  +4.96%  "
Set-TextValue 34 5 "  -3.39%  "
Set-TextValue 35 4 "41.37"
Set-TextValue 35 5 "  +0.75%  "
Set-TextValue 36 4 "1.00"
Set-TextValue 36 5 "  +0.01%  "
Set-TextValue 37 4 "56.77"
Set-TextValue 37 5 "  -2.82%  "
Set-TextValue 38 4 "0.0498"
Set-TextValue 38 5 "  +1.96%  "
Set-TextValue 39 4 "0.0₃0737"
Set-TextValue 39 5 "  +29.33%  "
Set-TextValue 40 4 "0.147"
Set-TextValue 40 5 "  +8.61%  "
Set-TextValue 41 4 "0.997"
Set-TextValue 41 5 "  -0.18%  "
Set-TextValue 42 4 "3.05"
Set-TextValue 42 5 "  -0.07%  "
Set-TextValue 43 2 "WEMIXToken"
Set-TextValue 43 3 "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue 43 4 "2.74"
Set-TextValue 43 5 "  -2.12%  "
Set-TextValue 44 2 "Monero"
Set-TextValue 44 3 "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue 44 4 "148.68"
Set-TextValue 44 5 "  -0.35%  "
Set-TextValue 45 4 "3.29"
Set-TextValue 45 5 "  -1.27%  "
Set-TextValue 46 4 "4.37"
Set-TextValue 46 5 "  -1.90%  "
Set-TextValue 47 5 "  -3.01%  "
Set-TextValue 48 4 "1.99"
Set-TextValue 48 5 "  -3.22%  "
Set-TextValue 49 4 "2.36"
Set-TextValue 49 5 "  -1.25%  "
Set-TextValue 50 4 "2.73"
Set-TextValue 50 5 "  +17.49%  "
Set-TextValue 51 4 "15.76"
Set-TextValue 51 5 "  -4.16%  "
